$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: new entry for "Shark" (in progress, not yet completed)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Shark"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "`n"
$ws.Range("G12").Value = 2019
$ws.Range("H12").Value = 12
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 13
$ws.Range("L12").Value = "None"
$ws.Range("M12").Value = "None"
$ws.Range("N12").Value = "None"
$ws.Range("O12").Value = "None"
$ws.Range("P12").Value = "None"
$ws.Range("Q12").Value = "None"

# Row 13: new entry for "Shark" (aborted)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Shark"
$ws.Range("C13").Value = 1126
$ws.Range("D13").Value = 1146
$ws.Range("E13").Value = 10
$ws.Range("F13").Value = "Very fast`n"
$ws.Range("G13").Value = 2019
$ws.Range("H13").Value = 12
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 14
$ws.Range("L13").Value = 2019
$ws.Range("M13").Value = 12
$ws.Range("N13").Value = 5
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 15
$ws.Range("Q13").Value = "aborted"

# Newline content triggers Excel's automatic row-height adjustment;
# restore the rows to their default (non-custom) height.
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
